$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5

# Row 8 updates
$ws.Range("G8").Value = 2.8
$ws.Range("I8").Value = 2.45
$ws.Range("L8").Value = 3.25
$ws.Range("X8").Value = 13
$ws.Range("AA8").Value = 23
$ws.Range("AB8").Value = 34
$ws.Range("AH8").Value = 12
$ws.Range("AP8").Value = 26
$ws.Range("AW8").Value = 4.5
